# Adding data dictionary ("Column dictionary") sheet to the workbook.
$wb = $excel.ActiveWorkbook

# --- 1. Create the new "Column dictionary" worksheet as the LAST sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Column dictionary"

# --- 2. Fill in the data dictionary table (A1:D21) ---
$rows = @(
    @("Column name", "Explanation", "Format", "Required?"),
    @("Type", "Type of creature. Currently two types of characters are supported: Martials (e.g. Fighter) and Blasters (e.g. throw a fireball)", "Martial/Blaster", "Yes"),
    @("Name", "Name of creature.", "Text", "Yes"),
    @("HP", "Health of creature. All creatures are assumed to start at max health. ", "Number", "Yes"),
    @("AC", "Creature's AC", "Number", "Yes"),
    @("str_save", "Bonus to save", "Number", "Yes"),
    @("dex_save", "Bonus to save", "Number", "Yes"),
    @("con_save", "Bonus to save", "Number", "Yes"),
    @("wis_save", "Bonus to save", "Number", "Yes"),
    @("cha_save", "Bonus to save", "Number", "Yes"),
    @("int_save", "Bonus to save", "Number", "Yes"),
    @("initiative_bonus", "Bonus to initiative", "Number", "Yes"),
    @("healer", "Can this character heal? Takes values of True or False.", "True/False", "Yes"),
    @("heal_amount", "Healing amount, only relevant if the character can heal. If healer=True but amount is missing, error is raised.", "e.g. 1d4+4", "No"),
    @("number_of_attacks", "Number of attacks Martial creature can make.", "Number", "Martial only"),
    @("attack_bonus", "Bonus to attack", "Number", "Martial only"),
    @("number_of_targets", "Number of targets Blaster's spell attacks on average.", "Number", "Blaster only"),
    @("spell_save_dc", "Spell save DC (Blaster)", "Number", "Blaster only"),
    @("targeted_save", "Save that targets need to make against Blaster's spell (e.g. dex for fireball)", "str/dex/con/wis/cha/int", "Blaster only"),
    @("saved_damage", "Damage target takes if it succeeds the saving throw: 0 for no damage, 0.5 for half-damage", "Number", "Blaster only"),
    @("attack_damage", "Damage of attack/spell", "e.g. 2d6+1d4+4", "Yes"),
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowValues = $rows[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

# --- 3. Column widths (best-fit approximation) ---
$ws.Columns.Item(1).ColumnWidth = 17.78
$ws.Columns.Item(2).ColumnWidth = 109
$ws.Columns.Item(3).ColumnWidth = 22.22
$ws.Columns.Item(4).ColumnWidth = 9.89

# --- 4. Page setup ---
$ws.PageSetup.Orientation = 1

# --- 5. Fix up the "Heroes" sheet view: drop the frozen/scrolled topLeftCell
#        and select the header row A1:T1 instead of a single cell ---
$wsHeroes = $wb.Worksheets.Item("Heroes")
$wsHeroes.Activate()
$wsHeroes.Range("A1:T1").Select() | Out-Null

# --- 6. Make sure "Monsters" is no longer the selected/active tab ---
$wsMonsters = $wb.Worksheets.Item("Monsters")
$wsMonsters.Activate()
$wsMonsters.Range("F7").Select() | Out-Null

# --- 7. Finally activate the new "Column dictionary" sheet and select C22,
#        matching where the author clicked after finishing data entry ---
$ws.Activate()
$ws.Range("C22").Select() | Out-Null
